# ApplicantInfo.xlsx - add "Valid Candidate" boolean results in column F
# and leave the final selection on G10, matching the author's manual review pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column F holds the boolean RegEx validation result for each applicant row.
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $false
$ws.Range("F4").Value = $false
$ws.Range("F5").Value = $true
$ws.Range("F6").Value = $false

# Leave the cursor/selection where the author last clicked before saving.
$ws.Range("G10").Select() | Out-Null
